# Auto-generated market data refresh for Marilith_Profits workbook
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per scheduled runner
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 800
$ws.Range("I40").Value = 800
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 800
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -625
$ws.Range("N40").ClearContents()
$ws.Range("H51").Value = 100000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 100000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 100000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -100968
$ws.Range("H55").Value = 216.66667
$ws.Range("I55").Value = 216.66667
$ws.Range("K55").Value = 216.66667
$ws.Range("M55").Value = -2.666670000000011
$ws.Range("H111").Value = 653.2857
$ws.Range("I111").Value = 634.8
$ws.Range("K111").Value = 1904.4
$ws.Range("M111").Value = 1162.6
$ws.Range("H113").Value = 6149.4
$ws.Range("I113").Value = 4700
$ws.Range("J113").Value = 7598.8
$ws.Range("K113").Value = 4700
$ws.Range("L113").Value = 7598.8
$ws.Range("M113").Value = -1446
$ws.Range("N113").Value = -14106.8
$ws.Range("H116").Value = 5002.5
$ws.Range("I116").Value = 5002.5
$ws.Range("K116").Value = 5002.5
$ws.Range("M116").Value = -1560.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 47
$ws.Range("I5").Value = 74
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 74
$ws.Range("L5").Value = 20
$ws.Range("M5").Value = 38
$ws.Range("N5").Value = -244
$ws.Range("H61").Value = 2233.2222
$ws.Range("I61").Value = 2137.375
$ws.Range("K61").Value = 2137.375
$ws.Range("M61").Value = -1925.375
$ws.Range("H93").Value = 35224
$ws.Range("J93").Value = 50448
$ws.Range("L93").Value = 50448
$ws.Range("N93").Value = -55440
$ws.Range("H136").Value = 2233.2222
$ws.Range("I136").Value = 2137.375
$ws.Range("K136").Value = 6412.125
$ws.Range("M136").Value = -3862.125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 47
$ws.Range("I4").Value = 74
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 74
$ws.Range("L4").Value = 20
$ws.Range("M4").Value = 41
$ws.Range("N4").Value = -250
$ws.Range("H122").Value = 200000
$ws.Range("J122").Value = 200000
$ws.Range("L122").Value = 200000
$ws.Range("N122").Value = -209800
$ws.Range("H134").Value = 6905.2856
$ws.Range("I134").Value = 6674.4443
$ws.Range("K134").Value = 20023.3329
$ws.Range("M134").Value = -17488.3329

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10000
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H7").Value = 220.35715
$ws.Range("I7").Value = 92.77778000000001
$ws.Range("J7").Value = 450
$ws.Range("K7").Value = 92.77778000000001
$ws.Range("L7").Value = 450
$ws.Range("M7").Value = 20.22221999999999
$ws.Range("N7").Value = -676
$ws.Range("H31").Value = 1800.25
$ws.Range("I31").Value = 1697.5
$ws.Range("K31").Value = 1697.5
$ws.Range("M31").Value = -1402.5
$ws.Range("H34").Value = 1800.25
$ws.Range("I34").Value = 1697.5
$ws.Range("K34").Value = 1697.5
$ws.Range("M34").Value = -1495.5
$ws.Range("H58").Value = 1616.1428
$ws.Range("I58").Value = 1647.3077
$ws.Range("J58").Value = 1589.1333
$ws.Range("K58").Value = 1647.3077
$ws.Range("L58").Value = 1589.1333
$ws.Range("M58").Value = -1444.3077
$ws.Range("N58").Value = -1995.1333
$ws.Range("H99").Value = 2004.8334
$ws.Range("I99").Value = 2004.8334
$ws.Range("K99").Value = 2004.8334
$ws.Range("M99").Value = -506.8334
$ws.Range("H126").Value = 2004.8334
$ws.Range("I126").Value = 2004.8334
$ws.Range("K126").Value = 6014.5002
$ws.Range("M126").Value = -3544.5002
$ws.Range("H132").Value = 3568.889
$ws.Range("I132").Value = 3642.375
$ws.Range("K132").Value = 10927.125
$ws.Range("M132").Value = -8397.125
$ws.Range("H134").Value = 4474.5713
$ws.Range("I134").Value = 4356
$ws.Range("J134").Value = 4487.0527
$ws.Range("K134").Value = 13068
$ws.Range("L134").Value = 13461.1581
$ws.Range("M134").Value = -10533
$ws.Range("N134").Value = -18531.1581
$ws.Range("H136").Value = 1616.1428
$ws.Range("I136").Value = 1647.3077
$ws.Range("J136").Value = 1589.1333
$ws.Range("K136").Value = 4941.9231
$ws.Range("L136").Value = 4767.3999
$ws.Range("M136").Value = -2391.9231
$ws.Range("N136").Value = -9867.3999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1366.421
$ws.Range("I4").Value = 1328.875
$ws.Range("K4").Value = 3986.625
$ws.Range("M4").Value = -3874.625
$ws.Range("H68").Value = 1410.8889
$ws.Range("I68").Value = 1118
$ws.Range("J68").Value = 1996.6666
$ws.Range("K68").Value = 3354
$ws.Range("L68").Value = 5989.9998
$ws.Range("M68").Value = -2543
$ws.Range("N68").Value = -7611.9998
$ws.Range("H71").Value = 1410.8889
$ws.Range("I71").Value = 1118
$ws.Range("J71").Value = 1996.6666
$ws.Range("K71").Value = 10062
$ws.Range("L71").Value = 17969.9994
$ws.Range("M71").Value = -6006
$ws.Range("N71").Value = -26081.9994
$ws.Range("H74").Value = 377496.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 377496.5
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 1132489.5
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -1134611.5
$ws.Range("H77").Value = 377496.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 377496.5
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 3397468.5
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -3408076.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H36").Value = 4059.8
$ws.Range("I36").Value = 5999.5
$ws.Range("J36").Value = 2766.6667
$ws.Range("K36").Value = 5999.5
$ws.Range("L36").Value = 2766.6667
$ws.Range("M36").Value = -5514.5
$ws.Range("N36").Value = -3736.6667

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H46").Value = 4641.6206
$ws.Range("I46").Value = 7810.636
$ws.Range("J46").Value = 2705
$ws.Range("K46").Value = 7810.636
$ws.Range("L46").Value = 2705
$ws.Range("M46").Value = -7622.636
$ws.Range("N46").Value = -3081
$ws.Range("H136").Value = 3311.1667
$ws.Range("I136").Value = 2926.5
$ws.Range("J136").Value = 3849.7
$ws.Range("K136").Value = 8779.5
$ws.Range("L136").Value = 11549.1
$ws.Range("M136").Value = -6229.5
$ws.Range("N136").Value = -16649.1

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4204.7144
$ws.Range("J62").Value = 2000
$ws.Range("L62").Value = 2000
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 4204.7144
$ws.Range("J65").Value = 2000
$ws.Range("L65").Value = 10000
$ws.Range("N65").Value = -16240
$ws.Range("H68").Value = 40316
$ws.Range("J68").Value = 40316
$ws.Range("L68").Value = 40316
$ws.Range("N68").Value = -41938
$ws.Range("H71").Value = 40316
$ws.Range("J71").Value = 40316
$ws.Range("L71").Value = 120948
$ws.Range("N71").Value = -129060
$ws.Range("H126").Value = 33926.938
$ws.Range("I126").Value = 33910.23
$ws.Range("K126").Value = 101730.69
$ws.Range("M126").Value = -99260.69
$ws.Range("H136").Value = 5665.222
$ws.Range("I136").Value = 5735.8184
$ws.Range("J136").Value = 5554.2856
$ws.Range("K136").Value = 17207.4552
$ws.Range("L136").Value = 16662.8568
$ws.Range("M136").Value = -14657.4552
$ws.Range("N136").Value = -21762.8568
